$wb = $excel.ActiveWorkbook

# Update the "Last Updated" timestamp on the Metadata sheet
$metaWs = $wb.Worksheets.Item("Metadata")
$metaWs.Range("A2").Value = "05 Nov 2025, 11:25 AM"

# Update the "1 Year" column (F) values on the Industry Analysis sheet
$indWs = $wb.Worksheets.Item("Industry Analysis")
$indWs.Range("F2").Value = 21.0016
$indWs.Range("F3").Value = -16.2396
$indWs.Range("F4").Value = 27.1317
$indWs.Range("F5").Value = -50.6494
$indWs.Range("F6").Value = 53.2813
$indWs.Range("F7").Value = -8.106199999999999
$indWs.Range("F8").Value = -9.552099999999999
$indWs.Range("F9").Value = 36.3756
$indWs.Range("F10").Value = -6.1314
$indWs.Range("F11").Value = 31.9081
$indWs.Range("F12").Value = -18.4955
$indWs.Range("F13").Value = 14.0155
$indWs.Range("F14").Value = -36.0718
$indWs.Range("F15").Value = -0.1622
$indWs.Range("F16").Value = 0.1459
$indWs.Range("F17").Value = -22.0012
$indWs.Range("F18").Value = 1.0561
$indWs.Range("F19").Value = -27.708
$indWs.Range("F20").Value = 47.7309
$indWs.Range("F21").Value = 12.0959
$indWs.Range("F22").Value = 95.1491
$indWs.Range("F23").Value = -50.2657
$indWs.Range("F24").Value = -13.3427
$indWs.Range("F25").Value = -9.9316
$indWs.Range("F26").Value = 5.8244
$indWs.Range("F27").Value = -32.7692
$indWs.Range("F28").Value = -24.8224
$indWs.Range("F29").Value = -18.4191
$indWs.Range("F30").Value = 25.8569
$indWs.Range("F31").Value = 58.4712
$indWs.Range("F32").Value = -3.3862
$indWs.Range("F33").Value = -6.3282
$indWs.Range("F34").Value = 27.7203
$indWs.Range("F35").Value = 4.4873
$indWs.Range("F36").Value = -4.9458
$indWs.Range("F37").Value = 3.6074
$indWs.Range("F38").Value = -23.3973
$indWs.Range("F39").Value = 8.7355
$indWs.Range("F40").Value = -5.8541
$indWs.Range("F41").Value = -8.3934
$indWs.Range("F42").Value = 20.3818
$indWs.Range("F43").Value = 14.3164
$indWs.Range("F44").Value = -12.6846
$indWs.Range("F45").Value = 28.4075
$indWs.Range("F46").Value = -1.1135
$indWs.Range("F47").Value = -37.1997
$indWs.Range("F48").Value = -29.8569
$indWs.Range("F49").Value = -27.5511
$indWs.Range("F50").Value = -49.7478
$indWs.Range("F51").Value = -51.8002
$indWs.Range("F52").Value = -38.5254
$indWs.Range("F53").Value = -12.4886
$indWs.Range("F54").Value = -5.0725
$indWs.Range("F55").Value = -17.7445
$indWs.Range("F56").Value = -26.636
$indWs.Range("F57").Value = -29.3361
$indWs.Range("F58").Value = -11.9574
$indWs.Range("F59").Value = -24.5687
$indWs.Range("F60").Value = -12.3
$indWs.Range("F61").Value = -10.9446
$indWs.Range("F62").Value = -17.1229
$indWs.Range("F63").Value = -9.5038
$indWs.Range("F64").Value = 54.2749
$indWs.Range("F65").Value = -43.4736
$indWs.Range("F66").Value = 13.2687
$indWs.Range("F67").Value = 12.7149
$indWs.Range("F68").Value = 24.8057
$indWs.Range("F69").Value = -17.0328
$indWs.Range("F70").Value = -6.8927
$indWs.Range("F71").Value = 13.6034
$indWs.Range("F72").Value = 3.9995
$indWs.Range("F73").Value = -16.226
$indWs.Range("F74").Value = -16.2448
$indWs.Range("F75").Value = 28.6924
$indWs.Range("F76").Value = 48.9752
